# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet's tab name reference (the sheet XML "name" attribute)
$ws.Name = "SCD0018"

# Update TC_ID column (B) for all data rows from "DGS-294" to the new TC id
$ws.Range("B2").Value = "SCD0018-002"
$ws.Range("B3").Value = "SCD0018-002"
$ws.Range("B4").Value = "SCD0018-002"

# Column B needs to widen to fit the new, longer TC_ID values
$ws.Columns("B:B").ColumnWidth = 11.6

# Adjust selection to match the authored state
$ws.Range("B5").Select()
